# PollsData.xlsx update — "update w/ opinion way and cluster17 polls (1/11)"
#
# The "omit" flag (column AF) that marked a poll-variant row as excluded
# moves down by one row within each duplicated-poll group (157->... etc.),
# and two brand-new poll rows (cluster17 + opinionway) are appended at the
# bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- shift the "omit" (AF) marker within the existing duplicate groups ---

# group around row 157/158: the omit flag moves from 158 to 157
$ws.Range("AF157").Value = 1
$ws.Range("AF158").ClearContents()

# group around row 161/162: the omit flag moves from 162 to 161
$ws.Range("AF161").Value = 1
$ws.Range("AF162").ClearContents()

# group around row 163/164: the omit flag moves from 164 to 163
$ws.Range("AF163").Value = 1
$ws.Range("AF164").ClearContents()

# --- append the two new poll rows (cluster17, opinionway) ---

$ws.Range("A165").Value = 64
$ws.Range("B165").Value = 2022
$ws.Range("C165").Value = 19
$ws.Range("D165").Value = 1
$ws.Range("E165").Value = 6
$ws.Range("F165").Value = "cluster17"
$ws.Range("G165").Value = "online"
$ws.Range("H165").Value = "partially"
$ws.Range("I165").Value = 0
$ws.Range("J165").Value = 2192
$ws.Range("K165").Value = 1.5
$ws.Range("L165").Value = 0.5
$ws.Range("M165").Value = 12.5
$ws.Range("N165").Value = 2
$ws.Range("O165").Value = 0.5
$ws.Range("P165").Value = 4.5
$ws.Range("Q165").Value = 2
$ws.Range("R165").Value = 23
$ws.Range("S165").Value = 14
$ws.Range("V165").Value = 1
$ws.Range("W165").Value = 2
$ws.Range("X165").Value = 15
$ws.Range("Y165").Value = 13.5
$ws.Range("Z165").Value = 1
$ws.Range("AB165").Value = 1.5
$ws.Range("AD165").Value = 5.5

$ws.Range("A166").Value = 65
$ws.Range("B166").Value = 2022
$ws.Range("C166").Value = 20
$ws.Range("D166").Value = 1
$ws.Range("E166").Value = 10
$ws.Range("F166").Value = "opinionway"
$ws.Range("G166").Value = "online"
$ws.Range("H166").Value = "partially"
$ws.Range("I166").Value = 1
$ws.Range("J166").Value = 997
$ws.Range("K166").Value = 2
$ws.Range("L166").Value = 1
$ws.Range("M166").Value = 9
$ws.Range("N166").Value = 2
$ws.Range("O166").Value = 1
$ws.Range("P166").Value = 8
$ws.Range("Q166").Value = 4
$ws.Range("R166").Value = 25
$ws.Range("S166").Value = 16
$ws.Range("V166").Value = 2
$ws.Range("W166").Value = 2
$ws.Range("X166").Value = 17
$ws.Range("Y166").Value = 12

# --- leave the selection where the author left it (one row lower than
#     before, since two rows were appended) ---
$ws.Range("J167").Select()
